$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.02181713207182
$ws.Range("C2").Value = 15.01622286022091
$ws.Range("D2").Value = 16.41869698022617

$ws.Range("B3").Value = 1.953885087150319
$ws.Range("C3").Value = 2.331802943878921
$ws.Range("D3").Value = 3.111477058831733

$ws.Range("B4").Value = 0.3748338957133264
$ws.Range("C4").Value = 0.4340093184667319
$ws.Range("D4").Value = 0.5955044643359204

$ws.Range("B5").Value = 80.7430097365543
$ws.Range("C5").Value = 81.71042264687431
$ws.Range("D5").Value = 82.81703297948127
